# Update the "Metadata" sheet:
#  - Date (row 8) moves forward a few minutes
#  - Jurisdiction (row 11) gets a value of "FRANCE" (was blank)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B8").Value = "2025-07-11T12:29:53+00:00"
$ws.Range("B11").Value = "FRANCE"
